$plain = @{
    2 = @{ 'D'='71.617.37'; 'E'='  -0.70%  ' }
    3 = @{ 'D'='3.861.16'; 'E'='  -2.38%  ' }
    4 = @{ 'E'='  +0.05%  ' }
    5 = @{ 'E'='  +1.15%  ' }
    6 = @{ 'E'='  +9.73%  ' }
    7 = @{ 'E'='  -1.43%  ' }
    8 = @{ 'E'='  -0.02%  ' }
    9 = @{ 'E'='  +1.87%  ' }
    10 = @{ 'E'='  +6.19%  ' }
    11 = @{ 'E'='  -0.62%  ' }
    12 = @{ 'E'='  +0.71%  ' }
    13 = @{ 'E'='  +3.32%  ' }
    14 = @{ 'D'='4.489.26'; 'E'='  -2.38%  ' }
    15 = @{ 'D'='3.880.97'; 'E'='  -2.43%  ' }
    16 = @{ 'E'='  +2.44%  ' }
    17 = @{ 'E'='  -1.06%  ' }
    18 = @{ 'E'='  -4.80%  ' }
    19 = @{ 'B'='TRON'; 'C'='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; 'E'='  -1.63%  ' }
    20 = @{ 'B'='WrappedBTC'; 'C'='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; 'D'='71.511.40'; 'E'='  -0.79%  ' }
    21 = @{ 'E'='  +0.60%  ' }
    22 = @{ 'E'='  -1.17%  ' }
    23 = @{ 'E'='  -2.20%  ' }
    24 = @{ 'E'='  -6.38%  ' }
    25 = @{ 'E'='  -3.72%  ' }
    26 = @{ 'E'='  -4.72%  ' }
    27 = @{ 'E'='  -5.46%  ' }
    28 = @{ 'E'='  +0.14%  ' }
    29 = @{ 'E'='  -4.80%  ' }
    30 = @{ 'E'='  -3.78%  ' }
    31 = @{ 'E'='  -3.33%  ' }
    32 = @{ 'B'='Cosmos'; 'C'='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; 'E'='  -0.16%  ' }
    33 = @{ 'B'='InjectiveProtocol'; 'C'='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; 'E'='  -3.67%  ' }
    34 = @{ 'E'='  -5.33%  ' }
    35 = @{ 'D'='0.0₃0964'; 'E'='  +13.75%  ' }
    36 = @{ 'E'='  -3.17%  ' }
    37 = @{ 'E'='  -10.55%  ' }
    38 = @{ 'E'='  -5.22%  ' }
    39 = @{ 'E'='  -0.01%  ' }
    40 = @{ 'B'='FirstDigitalUSD'; 'C'='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; 'E'='  +0.14%  ' }
    41 = @{ 'B'='ThetaToken'; 'C'='https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'; 'E'='  -0.40%  ' }
    42 = @{ 'B'='Kaspa'; 'C'='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; 'E'='  -2.68%  ' }
    43 = @{ 'E'='  +40.38%  ' }
    44 = @{ 'E'='  -4.06%  ' }
    45 = @{ 'E'='  -8.53%  ' }
    46 = @{ 'E'='  -3.76%  ' }
    47 = @{ 'E'='  -5.96%  ' }
    48 = @{ 'B'='ApeXProtocol'; 'C'='https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'; 'E'='  -0.73%  ' }
    49 = @{ 'B'='WEMIXToken'; 'C'='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; 'E'='  -17.10%  ' }
    50 = @{ 'D'='2.808.55'; 'E'='  +0.91%  ' }
    51 = @{ 'E'='  +3.25%  ' }
}

$textforce = @{
    5 = @{ 'D'='589.05' }
    6 = @{ 'D'='165.62' }
    7 = @{ 'D'='0.665' }
    9 = @{ 'D'='0.759' }
    10 = @{ 'D'='0.177' }
    11 = @{ 'D'='53.97' }
    12 = @{ 'D'='0.0000319' }
    13 = @{ 'D'='11.15' }
    16 = @{ 'D'='20.85' }
    17 = @{ 'D'='13.74' }
    18 = @{ 'D'='1.20' }
    19 = @{ 'D'='0.129' }
    21 = @{ 'D'='434.40' }
    22 = @{ 'D'='4.63' }
    23 = @{ 'D'='93.15' }
    24 = @{ 'D'='3.24' }
    25 = @{ 'D'='13.70' }
    26 = @{ 'D'='4.16' }
    27 = @{ 'D'='10.86' }
    28 = @{ 'D'='5.91' }
    29 = @{ 'D'='10.05' }
    30 = @{ 'D'='34.74' }
    31 = @{ 'D'='7.66' }
    32 = @{ 'D'='13.41' }
    33 = @{ 'D'='48.87' }
    34 = @{ 'D'='0.124' }
    36 = @{ 'D'='66.59' }
    37 = @{ 'D'='607.87' }
    38 = @{ 'D'='0.415' }
    40 = @{ 'D'='1.00' }
    41 = @{ 'D'='3.29' }
    42 = @{ 'D'='0.142' }
    43 = @{ 'D'='3.12' }
    44 = @{ 'D'='0.0465' }
    45 = @{ 'D'='10.12' }
    46 = @{ 'D'='0.143' }
    47 = @{ 'D'='2.58' }
    48 = @{ 'D'='3.32' }
    49 = @{ 'D'='2.80' }
    51 = @{ 'D'='0.000275' }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($rowKey in $plain.Keys) {
    $row = [int]$rowKey
    $cellMap = $plain[$rowKey]
    foreach ($col in $cellMap.Keys) {
        $ws.Range("$col$row").Value = $cellMap[$col]
    }
}

foreach ($rowKey in $textforce.Keys) {
    $row = [int]$rowKey
    $cellMap = $textforce[$rowKey]
    foreach ($col in $cellMap.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cellMap[$col]
        $cell.Style = "Normal"
    }
}
